$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 23:22"

# Update the row 5 (Galicia) figures
$ws.Range("B5").Value = 31043
$ws.Range("C5").Value = 13063
$ws.Range("D5").Value = 14832
$ws.Range("E5").Value = 3148
